$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.623.85'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '1.887.95'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2953'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06777'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("D10").Value = '1.888.08'
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.28'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07240'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '91.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6771'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("E15").Value = '  +3.53%  '
$ws.Range("D16").Value = '30.619.21'
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007983'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = '2.132.02'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.821'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '191.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +35.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.073'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.335'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.90%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.328'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09044'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05197'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7519'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.751'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01839'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.671'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.148'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9336'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4424'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.736'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.593'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1342'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05853'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.63%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.753'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.79%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.438'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3925'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.94%  '
